$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Sending cluster ECs -> Target cluster MuSCs, with refreshed TPM-based metrics
$ws.Range("A2").Value = "ECs"
$ws.Range("D2").Value = "MuSCs"

$ws.Range("G2").Value = 0.1207093333333333
$ws.Range("H2").Value = 0.362128
$ws.Range("I2").Value = 0.1354557203266995
$ws.Range("J2").Value = 0.1354557203266995
$ws.Range("M2").Value = 0.01533833333333333
$ws.Range("N2").Value = 0.046015
$ws.Range("O2").Value = 1
$ws.Range("P2").Value = 1
$ws.Range("Q2").Value = 0.001851479991111111
$ws.Range("R2").Value = 0.01666331992
$ws.Range("S2").Value = 0.1354557203266995
$ws.Range("T2").Value = 0.1354557203266995

# Row 3: MuSCs -> MuSCs, updated TPM-based metrics
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 0.7704256666666667
$ws.Range("H3").Value = 2.311277
$ws.Range("I3").Value = 0.8645442796733005
$ws.Range("J3").Value = 0.8645442796733006
$ws.Range("M3").Value = 0.046015
$ws.Range("O3").Value = 1
$ws.Range("P3").Value = 1
$ws.Range("Q3").Value = 0.01181704568388889
$ws.Range("R3").Value = 0.106353411155
$ws.Range("S3").Value = 0.8645442796733005
$ws.Range("T3").Value = 0.8645442796733006
